$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Trends Status" - update values
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 4

$ws1.Range("B3").Value = 13
$ws1.Range("C3").Value = 19
$ws1.Range("D3").Value = 13.5
$ws1.Range("E3").Value = 15.3

$ws1.Range("B4").Value = 61
$ws1.Range("C4").Value = 93
$ws1.Range("D4").Value = 63.5
$ws1.Range("E4").Value = 75

$ws1.Range("B5").Value = 15
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 15.6
$ws1.Range("E5").Value = 3.2

$ws1.Range("B6").Value = 6
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 6.2
$ws1.Range("E6").Value = 2.4

$ws1.Range("B7").Value = 210
$ws1.Range("C7").Value = 320

$ws1.Range("B8").Value = 451
$ws1.Range("C8").Value = 313

# ---------------------------------------------------------------------------
# Sheet 3: "Priority Status" - update values
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# Sheet 4: "Species qualification" - update text + values
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 757

$ws4.Range("B3").Value = 306
$ws4.Range("C3").Value = 96

$ws4.Range("B4").Value = 444
$ws4.Range("C4").Value = 124

# ---------------------------------------------------------------------------
# Sheet 5: "High Priority break-up" -> split into two sheets:
#   - existing sheet (sheetId 5) renamed to "Interannual update - High Pri"
#     and refilled with new "interannual" numbers (gains a "Trend Different"
#     row)
#   - a brand-new sheet (sheetId 6) "Major update - High Priority " that
#     keeps the original numbers that used to live on "High Priority
#     break-up"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("High Priority break-up")

# Create the new sheet right after the current last sheet, and populate it
# with the values the old "High Priority break-up" sheet used to hold.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMajor = $wb.Worksheets.Add($null, $lastSheet)
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"

$wsMajor.Range("A2").Value = "Trend New"
$wsMajor.Range("B2").Value = 14
$wsMajor.Range("C2").Value = 14.6
$wsMajor.Range("D2").Value = 14
$wsMajor.Range("E2").Value = 14.6

$wsMajor.Range("A3").Value = "IUCN"
$wsMajor.Range("B3").Value = 82
$wsMajor.Range("C3").Value = 85.40000000000001
$wsMajor.Range("D3").Value = 82
$wsMajor.Range("E3").Value = 85.40000000000001

$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108   # xlCenter

# Insert a new row for "Trend Different" on the original sheet, rename it,
# and refresh all of its values to the new interannual-update numbers.
$ws5.Rows.Item(3).Insert()

$ws5.Name = "Interannual update - High Pri"

$ws5.Range("A2").Value = "Trend New"
$ws5.Range("B2").Value = 8
$ws5.Range("C2").Value = 7.8
$ws5.Range("D2").Value = 8
$ws5.Range("E2").Value = 26.7

$ws5.Range("A3").Value = "Trend Different"
$ws5.Range("B3").Value = 7
$ws5.Range("C3").Value = 6.8
$ws5.Range("D3").Value = 1
$ws5.Range("E3").Value = 3.3

$ws5.Range("A4").Value = "IUCN"
$ws5.Range("B4").Value = 88
$ws5.Range("C4").Value = 85.40000000000001
$ws5.Range("D4").Value = 21
$ws5.Range("E4").Value = 70
